# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.375.07"
$ws.Range("E2").Value = "  -0.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.614.53"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "213.33"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.16%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("E9").Value = "  -0.79%  "

# Row 10
$ws.Range("E10").Value = "  +1.89%  "

# Row 11
$ws.Range("E11").Value = "  -0.85%  "

# Row 12
$ws.Range("D12").Value = "1.838.36"
$ws.Range("E12").Value = "  +0.16%  "

# Row 13
$ws.Range("D13").Value = "1.623.64"
$ws.Range("E13").Value = "  +0.89%  "

# Row 14
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("E15").Value = "  +0.34%  "

# Row 16
$ws.Range("D16").Value = "26.367.10"
$ws.Range("E16").Value = "  +0.06%  "

# Row 17
$ws.Range("D17").Value = "61.99"
$ws.Range("E17").Value = "  +1.55%  "

# Row 18
$ws.Range("E18").Value = "  +0.16%  "

# Row 19
$ws.Range("E19").Value = "  +0.00%  "

# Row 20
$ws.Range("D20").Value = "202.84"
$ws.Range("E20").Value = "  -2.53%  "

# Row 21
$ws.Range("D21").Value = "4.29"
$ws.Range("E21").Value = "  +0.44%  "

# Row 22
$ws.Range("D22").Value = "9.34"
$ws.Range("E22").Value = "  -0.45%  "

# Row 23
$ws.Range("E23").Value = "  +0.10%  "

# Row 24
$ws.Range("E24").Value = "  +5.40%  "

# Row 25
$ws.Range("D25").Value = "144.77"
$ws.Range("E25").Value = "  +1.70%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("E27").Value = "  -2.89%  "

# Row 28
$ws.Range("E28").Value = "  -0.62%  "

# Row 29 - force text so the trailing zero in "6.60" survives (would otherwise
# be read back as the number 6.6).
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.60"
$ws.Range("E29").Value = "  +1.70%  "

# Row 30
$ws.Range("E30").Value = "  +3.88%  "

# Row 31
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.18%  "

# Row 33 (row 32 unchanged)
$ws.Range("D33").Value = "2.96"
$ws.Range("E33").Value = "  -2.07%  "

# Row 34
$ws.Range("E34").Value = "  +3.12%  "

# Row 35
$ws.Range("E35").Value = "  +0.54%  "

# Row 36
$ws.Range("D36").Value = "1.162.34"
$ws.Range("E36").Value = "  +4.80%  "

# Row 37
$ws.Range("E37").Value = "  +1.80%  "

# Row 38
$ws.Range("E38").Value = "  +0.11%  "

# Row 39
$ws.Range("D39").Value = "0.793"
$ws.Range("E39").Value = "  +0.60%  "

# Row 40 - was ImmutableX, is now MXToken (rows 40/41 swapped rank order)
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  -0.47%  "

# Row 41 - was MXToken, is now ImmutableX
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "0.504"
$ws.Range("E41").Value = "  +0.79%  "

# Row 42
$ws.Range("E42").Value = "  +0.95%  "

# Row 43
$ws.Range("E43").Value = "  +2.57%  "

# Row 44
$ws.Range("D44").Value = "1.752.59"
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$ws.Range("D45").Value = "91.71"
$ws.Range("E45").Value = "  -1.43%  "

# Row 46
$ws.Range("E46").Value = "  -1.67%  "

# Row 47
$ws.Range("D47").Value = "54.37"
$ws.Range("E47").Value = "  +1.02%  "

# Row 48
$ws.Range("E48").Value = "  +0.32%  "

# Row 49 - was Mantle, is now BabyDogeCoin (rows 49/50 swapped rank order)
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0990"
$ws.Range("E49").Value = "  -6.43%  "

# Row 50 - was BabyDogeCoin, is now Mantle
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.407"
$ws.Range("E50").Value = "  -0.52%  "

# Row 51
$ws.Range("E51").Value = "  +0.10%  "
